$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.014.27"
$ws.Range("E2").Value = "  +0.19%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.547.89"
$ws.Range("E3").Value = "  +3.64%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.24"
$ws.Range("E5").Value = "  +1.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.15"
$ws.Range("E6").Value = "  +0.89%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.545.62"
$ws.Range("E7").Value = "  +3.54%  "

# Row 8
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("E9").Value = "  +4.43%  "

# Row 10
$ws.Range("E10").Value = "  +1.47%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.94"
$ws.Range("E11").Value = "  -1.55%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.413"
$ws.Range("E12").Value = "  +1.61%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.148.02"
$ws.Range("E13").Value = "  +3.68%  "

# Row 14
$ws.Range("E14").Value = "  +2.65%  "

# Row 15
$ws.Range("E15").Value = "  +1.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.547.75"
$ws.Range("E16").Value = "  +3.84%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.097.53"
$ws.Range("E17").Value = "  +0.43%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.116"
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.30"
$ws.Range("E19").Value = "  +8.66%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("E20").Value = "  +0.71%  "

# Row 21
$ws.Range("E21").Value = "  +1.59%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "430.07"
$ws.Range("E22").Value = "  +3.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.610"
$ws.Range("E23").Value = "  +5.32%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.10"
$ws.Range("E24").Value = "  +2.24%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.688.40"
$ws.Range("E25").Value = "  +3.75%  "

# Row 26
$ws.Range("E26").Value = "  -0.04%  "

# Row 27
$ws.Range("E27").Value = "  +6.80%  "

# Row 28
$ws.Range("E28").Value = "  +3.62%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.93"
$ws.Range("E29").Value = "  +0.95%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.04"
$ws.Range("E30").Value = "  -2.51%  "

# Row 32
$ws.Range("E32").Value = "  +0.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.52"
$ws.Range("E33").Value = "  +3.54%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.542.18"
$ws.Range("E34").Value = "  +3.63%  "

# Row 35
$ws.Range("E35").Value = "  -4.70%  "

# Row 36
$ws.Range("E36").Value = "  +0.09%  "

# Row 37
$ws.Range("E37").Value = "  +3.06%  "

# Row 38
$ws.Range("E38").Value = "  +4.17%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.60"
$ws.Range("E39").Value = "  +1.71%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.14%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "173.62"
$ws.Range("E41").Value = "  +2.94%  "

# Row 42
$ws.Range("E42").Value = "  -0.76%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.20"
$ws.Range("E43").Value = "  +2.90%  "

# Row 44
$ws.Range("E44").Value = "  +1.84%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.91"
$ws.Range("E45").Value = "  +0.73%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.03"
$ws.Range("E46").Value = "  +1.25%  "

# Row 47
$ws.Range("E47").Value = "  +1.72%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.46"
$ws.Range("E48").Value = "  -3.18%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.42"
$ws.Range("E49").Value = "  +13.39%  "

# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.11"
$ws.Range("E50").Value = "  +0.64%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.32"
$ws.Range("E51").Value = "  +2.28%  "
